$d = $word.ActiveDocument
$p = $d.Paragraphs.Item(102)
$r = $p.Range
Write-Host "Before text: [$($r.Text)]"

$xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="5768D874" w14:textId="12591621" w:rsidR="00B65A2B" w:rsidRPr="00B65A2B" w:rsidRDefault="00B65A2B" w:rsidP="00B65A2B"><w:pPr><w:spacing w:line="360" w:lineRule="auto"/><w:rPr><w:rFonts w:asciiTheme="majorBidi" w:hAnsiTheme="majorBidi" w:cstheme="majorBidi"/><w:b/><w:bCs/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r w:rsidRPr="00B65A2B"><w:rPr><w:rFonts w:asciiTheme="majorBidi" w:hAnsiTheme="majorBidi" w:cstheme="majorBidi"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve">GitHub Repo: </w:t></w:r><w:r><w:rPr><w:rFonts w:asciiTheme="majorBidi" w:hAnsiTheme="majorBidi" w:cstheme="majorBidi"/><w:b/><w:bCs/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>https://github.com/hamza01012/Lab</w:t></w:r><w:r><w:rPr><w:rFonts w:asciiTheme="majorBidi" w:hAnsiTheme="majorBidi" w:cstheme="majorBidi"/><w:b/><w:bCs/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>-</w:t></w:r><w:r><w:rPr><w:rFonts w:asciiTheme="majorBidi" w:hAnsiTheme="majorBidi" w:cstheme="majorBidi"/><w:b/><w:bCs/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>work/tree/main/Project______Environmental%20Classification%20for%20Plants</w:t></w:r></w:p>'
$res = $r.InsertXML($xml)
Write-Host "InsertXML result: $res"

$full = $d.Paragraphs.Item(102).Range
Write-Host "Full text after: [$($full.Text)]"
Write-Host "Doc paragraphs count: $($d.Paragraphs.Count)"

$p103 = $d.Paragraphs.Item(103)
Write-Host "Para103 text: [$($p103.Range.Text)] Start: $($p103.Range.Start) End: $($p103.Range.End)"
$p103.Range.Delete()
Write-Host "Doc paragraphs count after Range.Delete: $($d.Paragraphs.Count)"

$p103c = $d.Paragraphs.Item(103)
$p103c.Range.Delete()
Write-Host "Doc paragraphs count after range delete 2nd try: $($d.Paragraphs.Count)"

